# NV-22 Nguyễn Phúc Nam 7-2024 — "Đơn sale chính" sheet
# Append a new data row (row 3) below the existing header (row 1) and
# single data row (row 2), mirroring the row-2 pattern but with blank
# text columns and the numeric totals shown in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Đơn sale chính")

$row = 3

# Text columns are left blank (empty string) on the new row.
$ws.Cells.Item($row, 1).Value  = ""   # A - Tiền tố
$ws.Cells.Item($row, 3).Value  = ""   # C - Ngày thực hiện
$ws.Cells.Item($row, 4).Value  = ""   # D - Cơ sở
$ws.Cells.Item($row, 5).Value  = ""   # E - Khách hàng
$ws.Cells.Item($row, 6).Value  = ""   # F - Nguồn khách
$ws.Cells.Item($row, 7).Value  = ""   # G - Tên dịch vụ
$ws.Cells.Item($row, 8).Value  = ""   # H - Sale chính
$ws.Cells.Item($row, 10).Value = ""   # J - Sale phụ
$ws.Cells.Item($row, 17).Value = ""   # Q - Bác sĩ 1
$ws.Cells.Item($row, 18).Value = ""   # R - Bác sĩ 2
$ws.Cells.Item($row, 19).Value = ""   # S - Phụ phẫu 1
$ws.Cells.Item($row, 20).Value = ""   # T - Phụ phẫu 2

# Numeric columns carry the values from the diff.
$ws.Cells.Item($row, 2).Value  = 0          # B  - Mã dịch vụ
$ws.Cells.Item($row, 9).Value  = 3000000    # I  - Đơn giá gốc
$ws.Cells.Item($row, 11).Value = 0          # K  - Upsale
$ws.Cells.Item($row, 12).Value = 3000000    # L  - Đơn giá
$ws.Cells.Item($row, 13).Value = 3000000    # M  - Thanh toán lần đầu
$ws.Cells.Item($row, 14).Value = 0          # N  - Trả sau
$ws.Cells.Item($row, 15).Value = 3000000    # O  - Đã thanh toán
$ws.Cells.Item($row, 16).Value = 0          # P  - Dư nợ
